# Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025
#
# This workbook's sheet/TC_ID naming is being renumbered from the old
# "DGS-320" / "SCD0305" scheme to the new "SCD0019-003" / "SCD0019" scheme.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new TC numbering (SCD0305 -> SCD0019)
$ws.Name = "SCD0019"

# Update the TC_ID value in B2 (was "DGS-320", now "SCD0019-003")
$ws.Range("B2").Value = "SCD0019-003"

# Column B (TC_ID) is set to best-fit its content; refresh the auto width
# now that the TC_ID text is longer than before.
$ws.Columns.Item(2).AutoFit() | Out-Null

# Move / record the active selection on the sheet to B3 (as last edited by the author)
$ws.Range("B3").Select() | Out-Null
